# "유진증권" IFRS sheet - correct the financial figures: rows 2-6 were
# populated from the wrong source cells/units and get restated, with a few
# now-unused metric cells cleared; rows 7-9 lose their financial columns
# entirely (D:AJ cleared), keeping only the row label / period columns A:C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: restated figures ---
$ws.Range("D2").Value = 5678
$ws.Range("E2").Value = 149
$ws.Range("F2").Value = 149
$ws.Range("G2").Value = 119
$ws.Range("H2").Value = 64
$ws.Range("I2").Value = 64
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 46632
$ws.Range("L2").Value = 40934
$ws.Range("M2").Value = 5698
$ws.Range("N2").Value = 5698
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 5376
$ws.Range("Q2").Value = -2557
$ws.Range("R2").Value = -57
$ws.Range("S2").Value = 2436
$ws.Range("T2").Value = 51
$ws.Range("V2").Value = 6921
$ws.Range("W2").Value = 2.63
$ws.Range("X2").Value = 1.13
$ws.Range("Y2").Value = 1.22
$ws.Range("Z2").Value = 0.15
$ws.Range("AA2").Value = 718.36
$ws.Range("AB2").Value = 6.05
$ws.Range("AC2").Value = 88
$ws.Range("AD2").Value = 22.36
$ws.Range("AE2").Value = 5884
$ws.Range("AF2").Value = 0.33
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 96866418
$ws.Range("U2").ClearContents()

# --- Row 3: restated figures ---
$ws.Range("D3").Value = 7176
$ws.Range("E3").Value = 612
$ws.Range("F3").Value = 612
$ws.Range("G3").Value = 606
$ws.Range("H3").Value = 519
$ws.Range("I3").Value = 519
$ws.Range("K3").Value = 55843
$ws.Range("L3").Value = 49708
$ws.Range("M3").Value = 6135
$ws.Range("N3").Value = 6135
$ws.Range("P3").Value = 5376
$ws.Range("Q3").Value = -450
$ws.Range("R3").Value = 22
$ws.Range("S3").Value = 505
$ws.Range("T3").Value = 74
$ws.Range("V3").Value = 7242
$ws.Range("W3").Value = 8.529999999999999
$ws.Range("X3").Value = 7.24
$ws.Range("Y3").Value = 8.779999999999999
$ws.Range("Z3").Value = 1.01
$ws.Range("AA3").Value = 810.21
$ws.Range("AB3").Value = 14.18
$ws.Range("AC3").Value = 536
$ws.Range("AD3").Value = 5.02
$ws.Range("AE3").Value = 6335
$ws.Range("AF3").Value = 0.42
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 96866418
$ws.Range("O3").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("U3").ClearContents()

# --- Row 4: restated figures ---
$ws.Range("D4").Value = 7152
$ws.Range("E4").Value = 613
$ws.Range("F4").Value = 613
$ws.Range("G4").Value = 609
$ws.Range("H4").Value = 460
$ws.Range("I4").Value = 460
$ws.Range("K4").Value = 58908
$ws.Range("L4").Value = 52264
$ws.Range("M4").Value = 6644
$ws.Range("N4").Value = 6644
$ws.Range("P4").Value = 5376
$ws.Range("Q4").Value = 456
$ws.Range("R4").Value = -12
$ws.Range("S4").Value = -264
$ws.Range("T4").Value = 32
$ws.Range("V4").Value = 7386
$ws.Range("W4").Value = 8.57
$ws.Range("X4").Value = 6.43
$ws.Range("Y4").Value = 7.2
$ws.Range("Z4").Value = 0.8
$ws.Range("AA4").Value = 786.66
$ws.Range("AB4").Value = 23.64
$ws.Range("AC4").Value = 475
$ws.Range("AD4").Value = 5.13
$ws.Range("AE4").Value = 6861
$ws.Range("AF4").Value = 0.35
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 96866418
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("U4").ClearContents()

# --- Row 5: restated figures ---
$ws.Range("D5").Value = 7318
$ws.Range("E5").Value = 714
$ws.Range("F5").Value = 714
$ws.Range("G5").Value = 712
$ws.Range("H5").Value = 561
$ws.Range("I5").Value = 561
$ws.Range("K5").Value = 72099
$ws.Range("L5").Value = 65019
$ws.Range("M5").Value = 7080
$ws.Range("N5").Value = 7080
$ws.Range("P5").Value = 5376
$ws.Range("Q5").Value = -2699
$ws.Range("R5").Value = 6
$ws.Range("S5").Value = 3511
$ws.Range("T5").Value = 19
$ws.Range("V5").Value = 11121
$ws.Range("W5").Value = 9.76
$ws.Range("X5").Value = 7.67
$ws.Range("Y5").Value = 8.18
$ws.Range("Z5").Value = 0.86
$ws.Range("AA5").Value = 918.36
$ws.Range("AB5").Value = 31.75
$ws.Range("AC5").Value = 579
$ws.Range("AD5").Value = 6.03
$ws.Range("AE5").Value = 7311
$ws.Range("AF5").Value = 0.48
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 96866418
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()
$ws.Range("U5").ClearContents()

# --- Row 6: restated figures ---
$ws.Range("D6").Value = 8165
$ws.Range("E6").Value = 659
$ws.Range("F6").Value = 659
$ws.Range("G6").Value = 653
$ws.Range("H6").Value = 465
$ws.Range("I6").Value = 465
$ws.Range("K6").Value = 72503
$ws.Range("L6").Value = 64986
$ws.Range("M6").Value = 7517
$ws.Range("N6").Value = 7517
$ws.Range("P6").Value = 5376
$ws.Range("Q6").Value = 4492
$ws.Range("R6").Value = -38
$ws.Range("S6").Value = -2185
$ws.Range("T6").Value = 43
$ws.Range("V6").Value = 8916
$ws.Range("W6").Value = 8.07
$ws.Range("X6").Value = 5.69
$ws.Range("Y6").Value = 6.37
$ws.Range("Z6").Value = 0.64
$ws.Range("AA6").Value = 864.53
$ws.Range("AB6").Value = 39.88
$ws.Range("AC6").Value = 480
$ws.Range("AD6").Value = 4.95
$ws.Range("AE6").Value = 7762
$ws.Range("AF6").Value = 0.31
$ws.Range("AG6").Value = 60
$ws.Range("AH6").Value = 2.53
$ws.Range("AI6").Value = 12.5
$ws.Range("AJ6").Value = 96866418
$ws.Range("U6").ClearContents()

# --- Rows 7-9: financial figures removed, only A:C (rank/name/period) remain ---
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
